$wb = $excel.ActiveWorkbook

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 199
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1962.5143
$ws.Range("I129").Value = 1953.4286
$ws.Range("J129").Value = 1964.7858
$ws.Range("K129").Value = 5860.2858
$ws.Range("L129").Value = 5894.357400000001
$ws.Range("M129").Value = -860.2857999999997
$ws.Range("N129").Value = -15894.3574

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 16899.35
$ws.Range("I132").Value = 2407.132
$ws.Range("K132").Value = 7221.396000000001
$ws.Range("M132").Value = -4691.396000000001

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4399.2383
$ws.Range("I137").Value = 1433.375
$ws.Range("K137").Value = 4300.125
$ws.Range("M137").Value = -1750.125

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1387.0862
$ws.Range("I74").Value = 1345.238
$ws.Range("J74").Value = 1496.9375
$ws.Range("K74").Value = 1345.238
$ws.Range("L74").Value = 1496.9375
$ws.Range("M74").Value = -471.2380000000001
$ws.Range("N74").Value = -3244.9375

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1387.0862
$ws.Range("I77").Value = 1345.238
$ws.Range("J77").Value = 1496.9375
$ws.Range("K77").Value = 6726.190000000001
$ws.Range("L77").Value = 7484.6875
$ws.Range("M77").Value = -2358.190000000001
$ws.Range("N77").Value = -16220.6875

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1661.125
$ws.Range("I97").Value = 1319.5
$ws.Range("J97").Value = 2002.75
$ws.Range("K97").Value = 1319.5
$ws.Range("L97").Value = 2002.75
$ws.Range("M97").Value = -823.5
$ws.Range("N97").Value = -2994.75

# ARM row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1714.1794
$ws.Range("I122").Value = 1690.3611
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5071.0833
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2621.0833
$ws.Range("N122").Value = -10900

# ARM row 127
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

# ARM row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 48032
$ws.Range("J131").Value = 48032
$ws.Range("L131").Value = 48032
$ws.Range("N131").Value = -58112

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9616982
$ws.Range("I132").Value = 13889911
$ws.Range("J132").Value = 2893.25
$ws.Range("K132").Value = 41669733
$ws.Range("L132").Value = 8679.75
$ws.Range("M132").Value = -41667203
$ws.Range("N132").Value = -13739.75

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 942.9
$ws.Range("I94").Value = 928.625
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 928.625
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -477.625
$ws.Range("N94").Value = -1902

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1938.5897
$ws.Range("I107").Value = 1877.76
$ws.Range("J107").Value = 2047.2142
$ws.Range("K107").Value = 1877.76
$ws.Range("L107").Value = 2047.2142
$ws.Range("M107").Value = 42.24000000000001
$ws.Range("N107").Value = -5887.2142

# BSM row 139
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 40296
$ws.Range("J139").Value = 48192.75
$ws.Range("L139").Value = 48192.75
$ws.Range("N139").Value = -58472.75

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2436.25
$ws.Range("I31").Value = 1133.5853
$ws.Range("J31").Value = 3341.4915
$ws.Range("K31").Value = 1133.5853
$ws.Range("L31").Value = 3341.4915
$ws.Range("M31").Value = -838.5853
$ws.Range("N31").Value = -3931.4915

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2436.25
$ws.Range("I34").Value = 1133.5853
$ws.Range("J34").Value = 3341.4915
$ws.Range("K34").Value = 1133.5853
$ws.Range("L34").Value = 3341.4915
$ws.Range("M34").Value = -931.5853
$ws.Range("N34").Value = -3745.4915

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1579.6389
$ws.Range("I134").Value = 1263.0667
$ws.Range("J134").Value = 3162.5
$ws.Range("K134").Value = 3789.2001
$ws.Range("L134").Value = 9487.5
$ws.Range("M134").Value = -1254.2001
$ws.Range("N134").Value = -14557.5

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2560.2712
$ws.Range("I5").Value = 4127.815
$ws.Range("J5").Value = 1237.6562
$ws.Range("K5").Value = 12383.445
$ws.Range("L5").Value = 3712.9686
$ws.Range("M5").Value = -12271.445
$ws.Range("N5").Value = -3936.9686

# CUL row 104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3962.5
$ws.Range("J104").Value = 3962.5
$ws.Range("L104").Value = 11887.5
$ws.Range("N104").Value = -17129.5

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6549.4116
$ws.Range("I113").Value = 13106.25
$ws.Range("J113").Value = 721.1111
$ws.Range("K113").Value = 39318.75
$ws.Range("L113").Value = 2163.3333
$ws.Range("M113").Value = -37148.75
$ws.Range("N113").Value = -6503.3333

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 33338344
$ws.Range("I134").Value = 47623324
$ws.Range("J134").Value = 6721.6665
$ws.Range("K134").Value = 142869972
$ws.Range("L134").Value = 20164.9995
$ws.Range("M134").Value = -142864902
$ws.Range("N134").Value = -30304.9995

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2560.2712
$ws.Range("I135").Value = 4127.815
$ws.Range("J135").Value = 1237.6562
$ws.Range("K135").Value = 37150.335
$ws.Range("L135").Value = 11138.9058
$ws.Range("M135").Value = -34615.335
$ws.Range("N135").Value = -16208.9058

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1411.1111
$ws.Range("I122").Value = 1616.6666
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4849.9998
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2399.9998
$ws.Range("N122").Value = -7900

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1890.6212
$ws.Range("I132").Value = 1223.2683
$ws.Range("J132").Value = 2985.08
$ws.Range("K132").Value = 3669.8049
$ws.Range("L132").Value = 8955.24
$ws.Range("M132").Value = -1139.8049
$ws.Range("N132").Value = -14015.24

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2810
$ws.Range("I93").Value = 3499
$ws.Range("J93").Value = 2711.5715
$ws.Range("K93").Value = 3499
$ws.Range("L93").Value = 2711.5715
$ws.Range("M93").Value = -2251
$ws.Range("N93").Value = -5207.5715

# LTW row 111
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 44383
$ws.Range("J111").Value = 44383
$ws.Range("L111").Value = 44383
$ws.Range("N111").Value = -52563

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 113212.11
$ws.Range("I122").Value = 144772
$ws.Range("J122").Value = 2752.5
$ws.Range("K122").Value = 434316
$ws.Range("L122").Value = 8257.5
$ws.Range("M122").Value = -431866
$ws.Range("N122").Value = -13157.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2565.302
$ws.Range("I132").Value = 1759.1082
$ws.Range("K132").Value = 5277.3246
$ws.Range("M132").Value = -2747.3246

# WVR row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# WVR row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1583.4807
$ws.Range("I132").Value = 1134.1945
$ws.Range("K132").Value = 3402.5835
$ws.Range("M132").Value = -872.5835000000002

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 22704.412
$ws.Range("I136").Value = 29979.617
$ws.Range("J136").Value = 2091.3333
$ws.Range("K136").Value = 89938.851
$ws.Range("L136").Value = 6273.999899999999
$ws.Range("M136").Value = -87388.851
$ws.Range("N136").Value = -11373.9999

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 37276.332
$ws.Range("J139").Value = 32540
$ws.Range("L139").Value = 32540
$ws.Range("N139").Value = -42820
